$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 6-9 (old rows for MuSCs/D=ECs, MuSCs/D=Resolving-Mac, Resolving-Mac/D=ECs, Resolving-Mac/D=Resolving-Mac)
$ws.Range("A6:T9").EntireRow.Delete()

# Row 2: ECs -> Resolving-Mac
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("G2").Value = 0.8417533333333332
$ws.Range("H2").Value = 2.52526
$ws.Range("I2").Value = 0.01079423211523897
$ws.Range("J2").Value = 0.01079423211523897
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.888791333333333
$ws.Range("N2").Value = 5.666374
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1.589896400804444
$ws.Range("R2").Value = 14.30906760724
$ws.Range("S2").Value = 0.01079423211523897
$ws.Range("T2").Value = 0.01079423211523897

# Row 3: ECs -> FAPs (A), D stays Resolving-Mac
$ws.Range("A3").Value = "FAPs"
$ws.Range("G3").Value = 69.05064766666666
$ws.Range("H3").Value = 207.151943
$ws.Range("I3").Value = 0.885471656726338
$ws.Range("J3").Value = 0.8854716567263378
$ws.Range("M3").Value = 1.888791333333333
$ws.Range("N3").Value = 5.666374
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 130.4222648738535
$ws.Range("R3").Value = 1173.800383864682
$ws.Range("S3").Value = 0.885471656726338
$ws.Range("T3").Value = 0.8854716567263378

# Row 4: FAPs -> MuSCs (A), D: ECs -> Resolving-Mac
$ws.Range("A4").Value = "MuSCs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("G4").Value = 7.697976666666666
$ws.Range("H4").Value = 23.09393
$ws.Range("I4").Value = 0.09871507918910555
$ws.Range("J4").Value = 0.09871507918910553
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.888791333333333
$ws.Range("N4").Value = 5.666374
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 14.53987161220222
$ws.Range("R4").Value = 130.85884450982
$ws.Range("S4").Value = 0.09871507918910555
$ws.Range("T4").Value = 0.09871507918910553

# Row 5: FAPs -> Resolving-Mac (A), D stays Resolving-Mac
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 0.391393
$ws.Range("H5").Value = 1.174179
$ws.Range("I5").Value = 0.005019031969317685
$ws.Range("J5").Value = 0.005019031969317684
$ws.Range("M5").Value = 1.888791333333333
$ws.Range("N5").Value = 5.666374
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.7392597063273335
$ws.Range("R5").Value = 6.653337356946
$ws.Range("S5").Value = 0.005019031969317685
$ws.Range("T5").Value = 0.005019031969317684
